$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @(0.06980006979693865, 94.19484456346231, 9.548713539599117)
    3  = @(0.055947019143138,   53.77434063864826, 7.87796277773428)
    4  = @(0.05905563088508203, 37.03478084743531, 3.851574749998576)
    5  = @(0.05306189368001298, 22.1694707001653,  2.859374786431)
    6  = @(0.06989414375058972, 15.65405413311668, 1.886950878413403)
    7  = @(0.0432446907781123,  12.33067512306539, 1.075724056044919)
    8  = @(0.04555984643783421, 10.06491400444265, 1.008604408737339)
    9  = @(0.07148489812788107, 7.81632186327964,  0.8141481613972915)
    10 = @(0.05900969317750723, 6.80837440574668,  0.815493858845263)
    11 = @(0.03079368780602751, 5.725781214687653, 0.3681758332542525)
    12 = @(0.02777814510892506, 4.231376724044816, 0.5672285566331765)
    13 = @(0.06962246288800464, 4.106750383412441, 0.5881245397233539)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
